$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    76617.89384858565,
    72241.62606426013,
    69258.72732073106,
    66478.10933061691,
    64802.30256497666,
    64096.37880141272,
    62720.70952149436,
    65319.41100251421,
    73222.02298507578,
    78890.28804140349,
    82305.32934475451,
    83992.2204709472,
    81708.85766991691,
    83537.80851053994,
    83825.26655059501,
    82992.72293425066,
    80757.75910449831,
    79091.69357466768,
    82601.34948879662,
    81771.57327572181,
    79777.6218110089,
    77050.55496722442,
    74922.36287793868,
    71819.88758842958
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
